$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NATMI LR-pair output recomputed (Efna5-Epha1) "following Dr Hou advice":
# each row's ligand-/receptor-expressing-cell counts changed from 1 to 3,
# which re-derives the average/total expression, specificity and edge-weight
# columns. Values below are the new canonical numbers for rows 2-7.

$values = @{
    2 = @{ E=3; G=1.666083666666667;  H=4.998251;   I=0.6125276070882968; J=0.6125276070882968; K=3; M=4.913477;          N=14.740431; O=0.3201889893918886; P=0.3201889893918885; Q=8.186263776242333;  R=73.67637398618101; S=0.1961245954882335;  T=0.1961245954882335  }
    3 = @{ E=3; G=1.666083666666667;  H=4.998251;   I=0.6125276070882968; J=0.6125276070882968; K=3; M=3.864911333333334;  N=11.594734; O=0.2518587252793198; P=0.2518587252793198; Q=6.439265645581555;  R=57.953390810234;   S=0.1542704223196505;  T=0.1542704223196505  }
    4 = @{ E=3; G=1.666083666666667;  H=4.998251;   I=0.6125276070882968; J=0.6125276070882968; K=3; M=6.567164333333333;  N=19.701493; O=0.4279522853287917; P=0.4279522853287916; Q=10.94144523208255;  R=98.47300708874299; S=0.2621325892804128;  T=0.2621325892804128  }
    5 = @{ E=3; G=1.053930333333333;  H=3.161791;   I=0.3874723929117032; J=0.3874723929117031; K=3; M=4.913477;          N=14.740431; O=0.3201889893918886; P=0.3201889893918885; Q=5.178462452435667;  R=46.60616207192101; S=0.124064393903655;   T=0.124064393903655   }
    6 = @{ E=3; G=1.053930333333333;  H=3.161791;   I=0.3874723929117032; J=0.3874723929117031; K=3; M=3.864911333333334;  N=11.594734; O=0.2518587252793198; P=0.2518587252793198; Q=4.073347289843778;  R=36.660125608594;   S=0.09758830295966932; T=0.09758830295966929 }
    7 = @{ E=3; G=1.053930333333333;  H=3.161791;   I=0.3874723929117032; J=0.3874723929117031; K=3; M=6.567164333333333;  N=19.701493; O=0.4279522853287917; P=0.4279522853287916; Q=6.921333694884779;  R=62.292003253963;   S=0.1658196960483789;  T=0.1658196960483788  }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value = $rowVals[$col]
    }
}
